$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 ("isa_template"): template metadata block
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Version bump
$ws1.Range("B4").Value = "1.0.3"

# Description: drop "and observation unit"
$ws1.Range("B5").Value = "Template to describe the biological material according to MIAPPE v.1.1."

# Tags row: "Observation Unit" -> "plant material" (new FOODON based tag)
$ws1.Range("D13").Value = "plant material"

# Tags Term Accession Number row: MIAPPE_0069 -> FOODON_00004331
$ws1.Range("D14").Value = "http://purl.obolibrary.org/obo/FOODON_00004331"

# Tags Term Source REF row: add FOODON as the term source for the new tag
$ws1.Range("D15").Value = "FOODON"

# ---------------------------------------------------------------------------
# Sheet 2 ("biological_material"): annotation table
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$lo = $ws2.ListObjects.Item(1)

# Update the Unit ontology references (UO) to the new bioregistry URL scheme
$ws2.Range("W2").Value = "https://bioregistry.io/UO:0000008"
$ws2.Range("AS2").Value = "https://bioregistry.io/UO:0000008"

# Remove the nine trailing "Observation unit" related columns:
#   Characteristic [Observation unit type] / Term Source REF (MIAPPE:0071) / Term Accession Number (MIAPPE:0071)
#   Characteristic [External ID]           / Term Source REF (MIAPPE:0072) / Term Accession Number (MIAPPE:0072)
#   Characteristic [Spatial distribution]  / Term Source REF (MIAPPE:0073) / Term Accession Number (MIAPPE:0073)
# These sit right before the final "Output [Sample Name]" column, so we
# capture that column's value, drop everything from the end down through
# those nine columns, and then re-append "Output [Sample Name]" as the new
# last column.
$outputColName = "Output [Sample Name]"
$outputValue = $lo.ListColumns.Item($lo.ListColumns.Count).DataBodyRange.Value2

for ($i = 0; $i -lt 10; $i++) {
    $last = $lo.ListColumns.Count
    $lo.ListColumns.Item($last).Delete()
}

$newCol = $lo.ListColumns.Add()
$newCol.Range.Cells(1).Value = $outputColName
$newCol.DataBodyRange.Value = $outputValue
